$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$meta = $wb.Worksheets.Item(1)

# Remove the duplicate "Contact" / "No display for ContactDetail" row (row 11)
$meta.Rows.Item(11).Delete()

# Update Version
$meta.Range("B3").Value = "6.0.0"

# Update Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# Replace the old "Contact" row (now row 10, after the delete above) with Jurisdiction info
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet 2: "Elements" ---
$elements = $wb.Worksheets.Item(2)

# The root Extension element's Short/Definition should reflect this specific extension,
# not the generic "Extension"/"An Extension" text
$elements.Range("K2").Value = "Employee Union Id"
$elements.Range("L2").Value = "Identification for the labor union member"
